$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 282.7879796666667
$ws.Range("H2").Value = 848.363939
$ws.Range("I2").Value = 0.9674521741401267
$ws.Range("J2").Value = 0.9674521741401266
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.09934133333334
$ws.Range("N2").Value = 63.29802400000001
$ws.Range("O2").Value = 0.2917236204149438
$ws.Range("P2").Value = 0.2917236204149438
$ws.Range("Q2").Value = 5966.640107950727
$ws.Range("R2").Value = 53699.76097155654
$ws.Range("S2").Value = 0.2822286508184664
$ws.Range("T2").Value = 0.2822286508184664

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 282.7879796666667
$ws.Range("H3").Value = 848.363939
$ws.Range("I3").Value = 0.9674521741401267
$ws.Range("J3").Value = 0.9674521741401266
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.81943766666667
$ws.Range("N3").Value = 107.458313
$ws.Range("O3").Value = 0.4952465516465762
$ws.Range("P3").Value = 0.4952465516465762
$ws.Range("Q3").Value = 10129.30641055277
$ws.Range("R3").Value = 91163.7576949749
$ws.Range("S3").Value = 0.4791273531258807
$ws.Range("T3").Value = 0.4791273531258807

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 282.7879796666667
$ws.Range("H4").Value = 848.363939
$ws.Range("I4").Value = 0.9674521741401267
$ws.Range("J4").Value = 0.9674521741401266
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 15.40769666666667
$ws.Range("N4").Value = 46.22309
$ws.Range("O4").Value = 0.2130298279384801
$ws.Range("P4").Value = 0.2130298279384801
$ws.Range("Q4").Value = 4357.111411683502
$ws.Range("R4").Value = 39214.00270515151
$ws.Range("S4").Value = 0.2060961701957797
$ws.Range("T4").Value = 0.2060961701957797

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.714696666666668
$ws.Range("H5").Value = 23.14409
$ws.Range("I5").Value = 0.02639291836872237
$ws.Range("J5").Value = 0.02639291836872237
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.09934133333334
$ws.Range("N5").Value = 63.29802400000001
$ws.Range("O5").Value = 0.2917236204149438
$ws.Range("P5").Value = 0.2917236204149438
$ws.Range("Q5").Value = 162.7750182531289
$ws.Range("R5").Value = 1464.97516427816
$ws.Range("S5").Value = 0.007699437699839762
$ws.Range("T5").Value = 0.007699437699839762

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.714696666666668
$ws.Range("H6").Value = 23.14409
$ws.Range("I6").Value = 0.02639291836872237
$ws.Range("J6").Value = 0.02639291836872237
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.81943766666667
$ws.Range("N6").Value = 107.458313
$ws.Range("O6").Value = 0.4952465516465762
$ws.Range("P6").Value = 0.4952465516465762
$ws.Range("Q6").Value = 276.3360963689078
$ws.Range("R6").Value = 2487.02486732017
$ws.Range("S6").Value = 0.01307100180999933
$ws.Range("T6").Value = 0.01307100180999933

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.714696666666668
$ws.Range("H7").Value = 23.14409
$ws.Range("I7").Value = 0.02639291836872237
$ws.Range("J7").Value = 0.02639291836872237
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 15.40769666666667
$ws.Range("N7").Value = 46.22309
$ws.Range("O7").Value = 0.2130298279384801
$ws.Range("P7").Value = 0.2130298279384801
$ws.Range("Q7").Value = 118.8657061153445
$ws.Range("R7").Value = 1069.7913550381
$ws.Range("S7").Value = 0.005622478858883277
$ws.Range("T7").Value = 0.005622478858883277

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.799090333333333
$ws.Range("H8").Value = 5.397271
$ws.Range("I8").Value = 0.006154907491150983
$ws.Range("J8").Value = 0.006154907491150983
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.09934133333334
$ws.Range("N8").Value = 63.29802400000001
$ws.Range("O8").Value = 0.2917236204149438
$ws.Range("P8").Value = 0.2917236204149438
$ws.Range("Q8").Value = 37.95962103250044
$ws.Range("R8").Value = 341.636589292504
$ws.Range("S8").Value = 0.001795531896637623
$ws.Range("T8").Value = 0.001795531896637623

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.799090333333333
$ws.Range("H9").Value = 5.397271
$ws.Range("I9").Value = 0.006154907491150983
$ws.Range("J9").Value = 0.006154907491150983
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.81943766666667
$ws.Range("N9").Value = 107.458313
$ws.Range("O9").Value = 0.4952465516465762
$ws.Range("P9").Value = 0.4952465516465762
$ws.Range("Q9").Value = 64.44240405153589
$ws.Range("R9").Value = 579.9816364638231
$ws.Range("S9").Value = 0.003048196710696204
$ws.Range("T9").Value = 0.003048196710696204

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.799090333333333
$ws.Range("H10").Value = 5.397271
$ws.Range("I10").Value = 0.006154907491150983
$ws.Range("J10").Value = 0.006154907491150983
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 15.40769666666667
$ws.Range("N10").Value = 46.22309
$ws.Range("O10").Value = 0.2130298279384801
$ws.Range("P10").Value = 0.2130298279384801
$ws.Range("Q10").Value = 27.71983813193222
$ws.Range("R10").Value = 249.47854318739
$ws.Range("S10").Value = 0.001311178883817156
$ws.Range("T10").Value = 0.001311178883817156
